$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7204
$ws.Range("C3").Value = 176355
$ws.Range("C4").Value = 166311
$ws.Range("C8").Value = 64.61
